# "cambios de las fracciones" - update the reporting period (Q3 2022 -> Q4 2022)
# in the "Reporte de Formatos" sheet and refresh the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Fecha de inicio del periodo que se informa (B8): 2022-07-01 -> 2022-10-01
$ws.Range("B8").Value = 44835
# Fecha de término del periodo que se informa (C8): 2022-09-30 -> 2022-12-31
$ws.Range("C8").Value = 44926
# Fecha de validación (F8): 2022-10-10 -> 2023-01-10
$ws.Range("F8").Value = 44936
# Fecha de actualización (G8): 2022-10-10 -> 2023-01-10
$ws.Range("G8").Value = 44936

# Move the active selection to reflect where the editor was working.
$ws.Range("E2").Select()
$ws.Range("E12").Select()
